# "new excel for geography"
#
# The мountains reference sheet ("Данные") had a bogus/duplicate entry for
# Australia: "Водораздельный хребет" (Great Dividing Range) at row 4, with
# height/first-ascent year both wrongly set to 1854 (a data-entry mistake).
# Removing that row is the substance of the edit; everything else (shared
# string table shrinking, the now-unused "placeholder" font/style, the
# autofilter/_FilterDatabase range, and the selection left on the
# newly-shifted row 4) falls out of actually deleting the row the way a
# person would in the UI: select the whole row, then delete it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Данные")
$ws.Activate()

# Select row 4 (Водораздельный хребет) and delete it, shifting rows 5-17 up.
$ws.Rows.Item(4).Select()
$ws.Rows.Item(4).Delete()

# The autofilter range shrinks from A1:E16 to A1:E15 along with the data.
$ws.AutoFilterMode = $false
$ws.Range("A1:E15").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $ws.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Данные!`$A`$1:`$E`$15"
    }
}
